$wb = $excel.ActiveWorkbook

# --- Sheet "Mensal" (sheet1): append one monthly summary row (row 14) ---
$wsMensal = $wb.Worksheets.Item("Mensal")

$wsMensal.Cells.Item(14, 1).Value = 44150
$wsMensal.Cells.Item(13, 1).Copy()
$wsMensal.Cells.Item(14, 1).PasteSpecial(-4122)   # xlPasteFormats - match date cell style (s=2)

$wsMensal.Cells.Item(14, 2).Value = 218.9
$wsMensal.Cells.Item(14, 3).Value = 278.3
$wsMensal.Cells.Item(14, 4).Value = -21.35

# --- Sheet "Diario" (sheet2): append 15 daily rows (368-382) ---
$wsDiario = $wb.Worksheets.Item("Diario")

$diarioData = @(
    @(44136, 205.1,  278.3, -26.3),
    @(44137, 209.51, 278.3, -24.72),
    @(44138, 208.94, 278.3, -24.92),
    @(44139, 209,    278.3, -24.9),
    @(44140, 211.52, 278.3, -24),
    @(44141, 209.4,  278.3, -24.76),
    @(44142, 210.26, 278.3, -24.45),
    @(44143, 211.64, 278.3, -23.95),
    @(44144, 216.45, 278.3, -22.22),
    @(44145, 216.97, 278.3, -22.04),
    @(44146, 227.14, 278.3, -18.38),
    @(44147, 229.17, 278.3, -17.65),
    @(44148, 232.52, 278.3, -16.45),
    @(44149, 240.94, 278.3, -13.43),
    @(44150, 244.9,  278.3, -12)
)

$startRow = 368
for ($i = 0; $i -lt $diarioData.Count; $i++) {
    $r = $startRow + $i
    $row = $diarioData[$i]

    $wsDiario.Cells.Item($r, 1).Value = $row[0]
    $wsDiario.Cells.Item($r - 1, 1).Copy()
    $wsDiario.Cells.Item($r, 1).PasteSpecial(-4122)   # xlPasteFormats - match date cell style (s=2)

    $wsDiario.Cells.Item($r, 2).Value = $row[1]
    $wsDiario.Cells.Item($r, 3).Value = $row[2]
    $wsDiario.Cells.Item($r, 4).Value = $row[3]
}
